$d = $word.ActiveDocument

# Find empty paragraphs (no visible text, just the paragraph mark) whose
# style is BodyText or FirstParagraph, and remove them. These are the
# blank spacer paragraphs left over from the markdown conversion that the
# commit removes to tighten up spacing before the bookmarked headings.
$toDelete = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Style.NameLocal
    $text = $p.Range.Text
    $isEmpty = ($text -eq "`r") -or ($text -eq "`a`r")
    if ($isEmpty -and (($styleName -eq "Body Text") -or ($styleName -eq "First Paragraph"))) {
        $toDelete += $i
    }
}

# Delete from the last one to the first so earlier indices stay valid.
for ($j = $toDelete.Count - 1; $j -ge 0; $j--) {
    $idx = $toDelete[$j]
    $d.Paragraphs.Item($idx).Range.Delete()
}
